$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "MIN"
$ws.Range("B10").Formula = "=MIN(B4:B7)"

$ws.Range("A11").Value = "MAX"
$ws.Range("B11").Formula = "=MAX(B4:B7)"

$ws.Range("A12").Value = "AVERAGE"
$ws.Range("B12").Formula = "=AVERAGE(B4:B7)"

$ws.Range("A13").Value = "COUNT"
$ws.Range("B13").Formula = "=COUNT(B4:B7)"

$ws.Range("D12").Select()
